$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 6 table: switch the table style to the new built-in style GUID.
#    (PowerPoint table styles must be changed via Table.ApplyStyle, a direct
#    assignment to Table.Style throws "cannot be assigned through a
#    property".)
# ---------------------------------------------------------------------------
$tableSlide = $p.Slides.Item(6)
foreach ($shp in $tableSlide.Shapes) {
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{98864103-CA7D-4407-BAFC-F5B0C9654E3E}")
    }
}

# ---------------------------------------------------------------------------
# 2) Re-theme the deck: swap the custom "Integral" theme for the default
#    "Office Theme" palette (Design tab -> Office Theme). Theme colours are
#    exposed on the slide/master's ThemeColorScheme collection, ordered
#    dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink. COM RGB values are
#    packed 0x00BBGGRR (the classic OLE RGB() layout), so build each value
#    from its R/G/B bytes instead of pasting the hex string directly.
# ---------------------------------------------------------------------------
function OleColor($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$officeThemeColors = @(
    @(0x00, 0x00, 0x00),   # 1  dk1
    @(0xFF, 0xFF, 0xFF),   # 2  lt1
    @(0x44, 0x54, 0x6A),   # 3  dk2
    @(0xE7, 0xE6, 0xE6),   # 4  lt2
    @(0x5B, 0x9B, 0xD5),   # 5  accent1
    @(0xED, 0x7D, 0x31),   # 6  accent2
    @(0xA5, 0xA5, 0xA5),   # 7  accent3
    @(0xFF, 0xC0, 0x00),   # 8  accent4
    @(0x44, 0x72, 0xC4),   # 9  accent5
    @(0x70, 0xAD, 0x47),   # 10 accent6
    @(0x05, 0x63, 0xC1),   # 11 hlink
    @(0x95, 0x4F, 0x72)    # 12 folHlink
)

$master = $p.Slides.Item(1).Master
$colorScheme = $master.Theme.ThemeVariants
$themeColors = $p.Slides.Item(1).ThemeColorScheme

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $rgb = $officeThemeColors[$i - 1]
    $themeColors.Item($i).RGB = OleColor $rgb[0] $rgb[1] $rgb[2]
}

Write-Output "table style + theme colors updated"
